$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coin name / link (text) columns: plain assignment is safe, these are not
# numeric-looking strings.
# Price / volume (D, E) columns: prefix with a literal leading apostrophe so
# Excel stores the exact text (e.g. "41.50", "93.090.28", "  -2.89%  ") as a
# quote-prefixed Text cell instead of silently parsing it into a Number and
# losing formatting (trailing zeros, thousand-dot grouping, padding spaces).

$ws.Range("D2").Value = "'93.090.28"
$ws.Range("E2").Value = "'  -2.89%  "
$ws.Range("D3").Value = "'3.306.81"
$ws.Range("E3").Value = "'  -4.96%  "
$ws.Range("E4").Value = "'  +0.21%  "
$ws.Range("D5").Value = "'228.44"
$ws.Range("E5").Value = "'  -5.60%  "
$ws.Range("D6").Value = "'613.61"
$ws.Range("E6").Value = "'  -5.15%  "
$ws.Range("D7").Value = "'1.36"
$ws.Range("E7").Value = "'  -7.69%  "
$ws.Range("D8").Value = "'0.381"
$ws.Range("E8").Value = "'  -5.76%  "
$ws.Range("E9").Value = "'  +0.10%  "
$ws.Range("D10").Value = "'0.914"
$ws.Range("E10").Value = "'  -8.62%  "
$ws.Range("D11").Value = "'3.304.80"
$ws.Range("E11").Value = "'  -4.97%  "
$ws.Range("D12").Value = "'41.50"
$ws.Range("E12").Value = "'  -1.01%  "
$ws.Range("D13").Value = "'0.191"
$ws.Range("E13").Value = "'  -3.70%  "
$ws.Range("B14").Value = "Toncoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D14").Value = "'5.92"
$ws.Range("E14").Value = "'  -3.46%  "
$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").Value = "'93.090.18"
$ws.Range("E15").Value = "'  -2.45%  "
$ws.Range("D16").Value = "'3.929.02"
$ws.Range("E16").Value = "'  -5.07%  "
$ws.Range("D17").Value = "'0.0000241"
$ws.Range("E17").Value = "'  -7.07%  "
$ws.Range("D18").Value = "'7.94"
$ws.Range("E18").Value = "'  -6.85%  "
$ws.Range("D19").Value = "'3.317.30"
$ws.Range("E19").Value = "'  -4.51%  "
$ws.Range("D20").Value = "'17.03"
$ws.Range("E20").Value = "'  -5.67%  "
$ws.Range("D21").Value = "'10.75"
$ws.Range("E21").Value = "'  -7.90%  "
$ws.Range("D22").Value = "'3.41"
$ws.Range("E22").Value = "'  +6.44%  "
$ws.Range("D23").Value = "'488.62"
$ws.Range("E23").Value = "'  -3.61%  "
$ws.Range("D24").Value = "'0.437"
$ws.Range("E24").Value = "'  -15.56%  "
$ws.Range("D25").Value = "'0.0000179"
$ws.Range("E25").Value = "'  -7.46%  "
$ws.Range("D26").Value = "'5.99"
$ws.Range("E26").Value = "'  -9.94%  "
$ws.Range("D27").Value = "'89.20"
$ws.Range("E27").Value = "'  -6.54%  "
$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").Value = "'3.505.07"
$ws.Range("E28").Value = "'  -4.30%  "
$ws.Range("B29").Value = "Aptos"
$ws.Range("C29").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D29").Value = "'11.55"
$ws.Range("E29").Value = "'  -4.97%  "
$ws.Range("E30").Value = "'  -0.04%  "
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").Value = "'10.90"
$ws.Range("E31").Value = "'  -7.18%  "
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").Value = "'0.136"
$ws.Range("E32").Value = "'  -1.02%  "
$ws.Range("D33").Value = "'2.60"
$ws.Range("E33").Value = "'  -6.09%  "
$ws.Range("D34").Value = "'1.01"
$ws.Range("E34").Value = "'  +0.19%  "
$ws.Range("D35").Value = "'0.172"
$ws.Range("E35").Value = "'  -7.05%  "
$ws.Range("D36").Value = "'27.97"
$ws.Range("E36").Value = "'  -10.79%  "
$ws.Range("D37").Value = "'0.522"
$ws.Range("E37").Value = "'  -9.16%  "
$ws.Range("D38").Value = "'531.43"
$ws.Range("E38").Value = "'  +1.05%  "
$ws.Range("E39").Value = "'  -0.07%  "
$ws.Range("D40").Value = "'7.27"
$ws.Range("E40").Value = "'  -7.26%  "
$ws.Range("D41").Value = "'0.146"
$ws.Range("E41").Value = "'  -3.27%  "
$ws.Range("D42").Value = "'1.34"
$ws.Range("E42").Value = "'  -7.52%  "
$ws.Range("D43").Value = "'0.849"
$ws.Range("E43").Value = "'  -7.77%  "
$ws.Range("D44").Value = "'24.01"
$ws.Range("E44").Value = "'  -0.51%  "
$ws.Range("D45").Value = "'3.59"
$ws.Range("E45").Value = "'  +2.33%  "
$ws.Range("B46").Value = "ImmutableX"
$ws.Range("C46").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D46").Value = "'1.65"
$ws.Range("E46").Value = "'  -3.83%  "
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").Value = "'0.0407"
$ws.Range("E47").Value = "'  -2.68%  "
$ws.Range("D48").Value = "'5.29"
$ws.Range("E48").Value = "'  -6.50%  "
$ws.Range("B49").Value = "OKB"
$ws.Range("C49").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D49").Value = "'51.92"
$ws.Range("E49").Value = "'  -3.30%  "
$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").Value = "'2.07"
$ws.Range("E50").Value = "'  -4.12%  "
$ws.Range("D51").Value = "'7.83"
$ws.Range("E51").Value = "'  -3.42%  "
